# Update the marksheet for 1401MM07: marking-per-correct-answer changed
# from 3 to 5, which changes the total correct marks and the
# "correct/total" summary string on the Total row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 "Marking": marks awarded per correct answer
$ws.Range("B11").Value = 5

# Row 12 "Total": total marks obtained (Right * Marking)
$ws.Range("B12").Value = 125

# Row 12 "Max" column: correct/total marks summary text
$ws.Range("E12").Value = "125/140"
